$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: change Status from FAIL to PASS
$ws.Range("C3").Value = "PASS"

# Row 4: change from TC-TESTVALIDLOGIN2/testValidLogin2/FAIL to TC-TESTVALIDLOGIN3/testValidLogin3/PASS
$ws.Range("A4").Value = "TC-TESTVALIDLOGIN3"
$ws.Range("B4").Value = "testValidLogin3"
$ws.Range("C4").Value = "PASS"

# New row 5: TC-TESTVALIDLOGIN3/testValidLogin3/PASS
$ws.Range("A5").Value = "TC-TESTVALIDLOGIN3"
$ws.Range("B5").Value = "testValidLogin3"
$ws.Range("C5").Value = "PASS"

# New row 6: TC-TESTVALIDLOGIN3/testValidLogin3/PASS
$ws.Range("A6").Value = "TC-TESTVALIDLOGIN3"
$ws.Range("B6").Value = "testValidLogin3"
$ws.Range("C6").Value = "PASS"
